$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# The new item "ماكينه حلاقه جليت فليكتور" (Gillette Flicker razor) needs to
# be inserted alphabetically as row 23 (item #17), pushing the existing
# "محلول ملح" row down to row 24 (item #18), and the totals / footer rows
# down from 24/25 to 25/26.
# -------------------------------------------------------------------------

# 1) Duplicate row 23 ("محلول ملح") into a new row 24 - this shifts the old
#    row 24 (totals) and row 25 (footer) down to 25 / 26, and gives the new
#    row the same styling/merge layout as the row it was cloned from.
$ws.Rows.Item(23).Copy()
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(24).RowHeight = 25.5

# Re-create the merged cells for the newly inserted row 24 (Insert does not
# carry merges over automatically).
$ws.Range("A24:B24").Merge()
$ws.Range("C24:G24").Merge()
$ws.Range("H24:K24").Merge()
$ws.Range("L24:M24").Merge()
$ws.Range("N24:O24").Merge()

# Row 24 is now an exact duplicate of (old) row 23 - just renumber it to 18.
$ws.Range("A24").Value = 18

# 2) Overwrite row 23 with the new product's data.
$ws.Range("C23").Value = "ماكينه حلاقه جليت فليكتور"
$ws.Range("H23").Value = "16:0"
$ws.Range("N23").Value = "15.00"
$ws.Range("Q23").Value = "2:0"
# L23 (order limit) stays "0", same as before - no change needed.

# P23 uses a numeric-looking number format ("#.00"), so a plain Value
# assignment would be auto-coerced to a number and lose the trailing
# zeros ("30.0000" -> 30). Flip to text, assign, then restore the
# original number format so the stored style id is unchanged.
$p23 = $ws.Range("P23")
$p23fmt = $p23.NumberFormat
$p23.NumberFormat = "@"
$p23.Value = "30.0000"
$p23.NumberFormat = $p23fmt

# 3) Update the running total in (now) row 25 to include the new item's
#    sell price (612.31 = 582.31 + 30.00).
$ws.Range("P25").Value = 612.31

# 4) Update the generated timestamp string (row 26 / A26 after the shift).
$ws.Range("A26").Value = "Saturday, 23 August, 2025 11:41 AM"

# Row 25 (totals row) height changes slightly once the new row is in place.
$ws.Rows.Item(25).RowHeight = 24.75
